# Update the author attribution line on the cheat sheet footer:
# add "Tiffany Timbers (https://www.tiffanytimbers.com/)" after the
# existing "Tomas Beuzen (https://www.tomasbeuzen.com/)" credit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the footer text box by its well-known shape Id (falls back to a
# name/content scan if the Id ever changes).
$targetId = 152
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq $targetId) {
        $shape = $candidate
        break
    }
}
if ($shape -eq $null) {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $candidate = $s.Shapes.Item($i)
        if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
            if ($candidate.TextFrame.TextRange.Text.Contains("tomasbeuzen.com")) {
                $shape = $candidate
                break
            }
        }
    }
}

$tr = $shape.TextFrame.TextRange
$fullText = $tr.Text

$existingUrl = "https://www.tomasbeuzen.com/"
$urlStart = $fullText.IndexOf($existingUrl)
$insertPos = $urlStart + $existingUrl.Length + 1

# The single character right after the existing hyperlinked URL is the
# closing ")" - insert the new attribution text immediately before it.
$closeParen = $tr.Characters($insertPos, 1)

$newUrl = "https://www.tiffanytimbers.com/"
$insertedText = ") & Tiffany Timbers (" + $newUrl
$closeParen.InsertBefore($insertedText) | Out-Null

# Re-fetch the text range/content now that the insertion shifted offsets,
# then apply the hyperlink to just the newly-added URL text.
$tr2 = $shape.TextFrame.TextRange
$fullText2 = $tr2.Text
$urlStart2 = $fullText2.IndexOf($newUrl)
$urlRange = $tr2.Characters($urlStart2 + 1, $newUrl.Length)
$urlRange.ActionSettings.Item(1).Hyperlink.Address = $newUrl
